$d = $word.ActiveDocument

function Find-ParagraphByExactText($doc, $searchText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $searchText) {
            return $i
        }
    }
    return -1
}

function Find-ParagraphStartingWith($doc, $prefix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

function Find-ParagraphEndingWith($doc, $suffix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t.EndsWith($suffix)) {
            return $i
        }
    }
    return -1
}

$pkgHeader = '<?xml version="1.0" standalone="yes"?>' + "`n" + '<?mso-application progid="Word.Document"?>' + "`n"
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1) Split-point / lastRenderedPageBreak shuffle inside the "bundling" para:
#    "Si bien por defecto la meta de este paso..." - the overall text is
#    unchanged, only where the run boundary (and lastRenderedPageBreak) sit.
# ---------------------------------------------------------------------------
$idx1 = Find-ParagraphStartingWith $d "Si bien por defecto la meta de este paso"
if ($idx1 -gt 0) {
    $p1 = $d.Paragraphs.Item($idx1)
    $full1 = $p1.Range
    $content1 = $d.Range($full1.Start, $full1.End - 1)

    $t1a = "Si bien por defecto la meta de este paso es agrupar todo en 3 archivos resultantes (uno de JS, CSS y HTML), puede darse el caso de que alguno de los archivos "
    $t1b = "resultante sea muy pesado, se considera pesado si pasa los 500 kb. Aquí es cuando es recomendable hacer un Split de la estrategia del bundling del empaquetador. Esto se puede lograr en la configuración del bundlers para declarar que paquetes o componentes deseamos que estén en un archivo diferente. Así logramos tener archivos mas pequeños, lo ideal es que estos adicionales no sean de uso primario, para que solo se carguen en el momento de su uso."

    $p1xml = '<w:p><w:r><w:rPr><w:lang w:val="es-EC"/></w:rPr><w:t xml:space="preserve">' + $t1a + '</w:t></w:r><w:r><w:rPr><w:lang w:val="es-EC"/></w:rPr><w:lastRenderedPageBreak/><w:t>' + $t1b + '</w:t></w:r></w:p>'
    $xml1 = $pkgHeader + $pkgOpen + $p1xml + $pkgClose
    [void]$content1.InsertXML($xml1)
}

# ---------------------------------------------------------------------------
# 2) "VS Code:" paragraph gains a <w:lastRenderedPageBreak/> before its text.
# ---------------------------------------------------------------------------
$idx2 = Find-ParagraphByExactText $d "VS Code:"
if ($idx2 -gt 0) {
    $p2 = $d.Paragraphs.Item($idx2)
    $full2 = $p2.Range
    $content2 = $d.Range($full2.Start, $full2.End - 1)

    $p2xml = '<w:p><w:r><w:rPr><w:lang w:val="es-EC"/></w:rPr><w:lastRenderedPageBreak/><w:t>VS Code:</w:t></w:r></w:p>'
    $xml2 = $pkgHeader + $pkgOpen + $p2xml + $pkgClose
    [void]$content2.InsertXML($xml2)
}

# ---------------------------------------------------------------------------
# 3) "Renderizado condicional y dinamico (Extra)" Heading1 loses the
#    <w:lastRenderedPageBreak/> that used to sit before its text.
# ---------------------------------------------------------------------------
$idx3 = Find-ParagraphStartingWith $d "Renderizado condicional y din"
if ($idx3 -gt 0) {
    $p3 = $d.Paragraphs.Item($idx3)
    $full3 = $p3.Range
    $content3 = $d.Range($full3.Start, $full3.End - 1)

    $p3xml = '<w:p><w:r><w:rPr><w:lang w:val="es-EC"/></w:rPr><w:t>Renderizado condicional y dinámico</w:t></w:r><w:r w:rsidR="0018041C"><w:rPr><w:lang w:val="es-EC"/></w:rPr><w:t xml:space="preserve"> (Extra)</w:t></w:r></w:p>'
    $xml3 = $pkgHeader + $pkgOpen + $p3xml + $pkgClose
    [void]$content3.InsertXML($xml3)
}

# ---------------------------------------------------------------------------
# 4) Insert a brand-new paragraph right before the "Compilacion" Heading1,
#    after the last "Props" paragraph ("... otro Componente de ReactJS.").
# ---------------------------------------------------------------------------
$idx4 = Find-ParagraphEndingWith $d "o incluso, otro Componente de ReactJS."
if ($idx4 -gt 0) {
    $p4 = $d.Paragraphs.Item($idx4)
    [void]$p4.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($idx4 + 1)
    $rng4 = $newPara.Range
    $rng4.Collapse(1)

    $text4 = "Como tal React no es difícil, pero si hay que tener en el mapa ciertas bases o conocimientos para construir buenas aplicaciones. Tales como el ciclo de vida y como utilizarlo, sobre todo el tema de memorization. Patrones de diseño eficientes. El uso de librerías adecuadas etc"

    $p4xml = '<w:p><w:pPr><w:rPr><w:lang w:val="es-EC"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-EC"/></w:rPr><w:t>' + $text4 + '</w:t></w:r><w:r><w:rPr><w:lang w:val="es-EC"/></w:rPr><w:t>.</w:t></w:r></w:p>'
    $xml4 = $pkgHeader + $pkgOpen + $p4xml + $pkgClose
    [void]$rng4.InsertXML($xml4)
}

Write-Host "Edits applied."
